# Weekly update: a new price-report row for "Femacal de La Calera" (Arveja
# Verde) is inserted at row 5, pushing all existing data rows (old 5..47)
# down by one (new 6..48). The new row 5 carries this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 5, shifting rows 5-47 -> 6-48.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with this period's data.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44550
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 100112022
$ws.Range("G5").Value = "Arveja Verde"
$ws.Range("H5").Value = "Perfection"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 73
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17521
$ws.Range("N5").Value = "`$/saco 25 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 701
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
